$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @($null, "id", "particle", "energy", "position", "direction", "event"),
  @(0, 1, "photon", 1000, "[1. 2. 3.]", "[15. 20.]", "PRIMARY"),
  @(1, 1, "photon", 1000, "[24508.65269744 54829.56344934 37690.39192712]", "[15. 20.]", "MOVE"),
  @(2, 2, "electron", 34.19096534435106, "[24508.65269744 54829.56344934 37690.39192712]", "[14.99948899  5.77235136]", "PAIR PRODUCTION"),
  @(3, 3, "positron", 965.8090346556489, "[24508.65269744 54829.56344934 37690.39192712]", "[14.99948899 -5.77235136]", "PAIR PRODUCTION"),
  @(4, 2, "electron", 31.89638555785016, "[24721.6900838  54710.16577208 38017.80126594]", "[0.32081834 1.19405907]", "MOVE"),
  @(5, 2, "electron", 17.5340151332472, "[24721.6900838  54710.16577208 38017.80126594]", "[0.32081834 1.19405907]", "BREMSSTRAHLUNG"),
  @(6, 4, "photon", 14.36237042460295, "[24721.6900838  54710.16577208 38017.80126594]", "[0.3368392  2.26060968]", "BREMSSTRAHLUNG"),
  @(7, 4, "photon", 14.36237042460295, "[19564.8616914  60960.697398   22414.95269621]", "[0.3368392  2.26060968]", "MOVE"),
  @(8, 5, "electron", 5.805174018462141, "[19564.8616914  60960.697398   22414.95269621]", "[0.37241882 3.65593015]", "PAIR PRODUCTION"),
  @(9, 6, "positron", 8.557196406140813, "[19564.8616914  60960.697398   22414.95269621]", "[ 0.37241882 -3.65593015]", "PAIR PRODUCTION"),
  @(10, 2, "electron", 5.770566295008839, "[24409.56940619 54367.19571873 37858.44258768]", "[0.33332945 2.54382368]", "MOVE"),
  @(11, 2, "electron", 2.987384164531072, "[24409.56940619 54367.19571873 37858.44258768]", "[0.33332945 2.54382368]", "BREMSSTRAHLUNG"),
  @(12, 7, "photon", 2.783182130477767, "[24409.56940619 54367.19571873 37858.44258768]", "[0.42188361 1.54811421]", "BREMSSTRAHLUNG"),
  @(13, 7, "photon", 2.783182130477767, "[  68061.00739401 1978524.49768041  144460.66410187]", "[0.42188361 1.54811421]", "MOVE"),
  @(14, 8, "electron", 1.724253354036095, "[  68061.00739401 1978524.49768041  144460.66410187]", "[0.23827808 0.30056747]", "PAIR PRODUCTION"),
  @(15, 9, "positron", 1.058928776441672, "[  68061.00739401 1978524.49768041  144460.66410187]", "[ 0.23827808 -0.30056747]", "PAIR PRODUCTION"),
  @(16, 3, "positron", 961.5091125054774, "[24879.50262682 55037.40781384 38260.33757025]", "[1.13217414 1.36691641]", "MOVE"),
  @(17, 3, "positron", 274.1870826078945, "[24879.50262682 55037.40781384 38260.33757025]", "[1.13217414 1.36691641]", "BREMSSTRAHLUNG"),
  @(18, 10, "photon", 687.322029897583, "[24879.50262682 55037.40781384 38260.33757025]", "[1.13164268 2.53692433]", "BREMSSTRAHLUNG"),
  @(19, 3, "positron", 266.7476377361234, "[25091.74805529 56063.97503394 38494.77543305]", "[0.29437794 1.24258164]", "MOVE"),
  @(20, 3, "positron", 245.5616149816444, "[25091.74805529 56063.97503394 38494.77543305]", "[0.29437794 1.24258164]", "BREMSSTRAHLUNG"),
  @(21, 11, "photon", 21.18602275447892, "[25091.74805529 56063.97503394 38494.77543305]", "[0.29246225 5.86720934]", "BREMSSTRAHLUNG"),
  @(22, 11, "photon", 21.18602275447892, "[29951.9990088  53916.93661228 55352.45522456]", "[0.29246225 5.86720934]", "MOVE"),
  @(23, 12, "electron", 16.56019496691658, "[29951.9990088  53916.93661228 55352.45522456]", "[0.31658228 3.63343614]", "PAIR PRODUCTION"),
  @(24, 13, "positron", 4.625827787562347, "[29951.9990088  53916.93661228 55352.45522456]", "[ 0.31658228 -3.63343614]", "PAIR PRODUCTION"),
  @(25, 12, "electron", 8.191997847693441, "[29882.1747087  54195.88760318 54772.39534386]", "[1.38707495 2.00859345]", "MOVE"),
  @(26, 12, "electron", 6.171394761063226, "[29882.1747087  54195.88760318 54772.39534386]", "[1.38707495 2.00859345]", "BREMSSTRAHLUNG"),
  @(27, 14, "photon", 2.020603086630215, "[29882.1747087  54195.88760318 54772.39534386]", "[1.44945382 0.87826284]", "BREMSSTRAHLUNG"),
  @(28, 14, "photon", 2.020603086630215, "[5855021.16886212 7075783.29386117 5923060.75478291]", "[1.44945382 0.87826284]", "MOVE"),
  @(29, 16, "positron", 1.094707468182343, "[5855021.16886212 7075783.29386117 5923060.75478291]", "[ 1.70235239 -3.59269334]", "PAIR PRODUCTION"),
  @(30, 3, "positron", 235.8105346386953, "[25234.04265409 56481.83558126 38985.20191172]", "[0.61465117 5.85135185]", "MOVE"),
  @(31, 3, "positron", 67.34670428513886, "[25234.04265409 56481.83558126 38985.20191172]", "[0.61465117 5.85135185]", "BREMSSTRAHLUNG"),
  @(32, 17, "photon", 168.4638303535564, "[25234.04265409 56481.83558126 38985.20191172]", "[0.6168182  0.98620691]", "BREMSSTRAHLUNG"),
  @(33, 3, "positron", 55.36382446076215, "[26247.57270091 56014.75837767 40742.74741594]", "[1.02290727 1.39072078]", "MOVE"),
  @(34, 18, "photon", 54.4624139141232, "[26247.57270091 56014.75837767 40742.74741594]", "[1.03213726 0.09969227]", "BREMSSTRAHLUNG"),
  @(35, 18, "photon", 54.4624139141232, "[58966.29564398 59287.41113592 78858.80599506]", "[1.03213726 0.09969227]", "MOVE"),
  @(36, 19, "electron", 44.25783988311611, "[58966.29564398 59287.41113592 78858.80599506]", "[1.04152001 5.20777295]", "PAIR PRODUCTION"),
  @(37, 20, "positron", 10.20457403100708, "[58966.29564398 59287.41113592 78858.80599506]", "[ 1.04152001 -5.20777295]", "PAIR PRODUCTION")
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $val = $row[$c]
    if ($null -ne $val) {
      $ws.Cells.Item($r + 1, $c + 1).Value = $val
    }
  }
}

Write-Host "done"